$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.385.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.839.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.04"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.21"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.833.03"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.45%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.70"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.45%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.477.92"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.853.97"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.595.62"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.24"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "483.11"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.716"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.83"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.90"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.988.46"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.36"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.94"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.783.50"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.85"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.317"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.95"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "429.48"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.35"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "142.97"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.824.63"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +13.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0355"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.56%  "
